# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (column F) counts on each sheet.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value2  = 7893
$ws.Cells.Item(5, 6).Value2  = 15241
$ws.Cells.Item(8, 6).Value2  = 654
$ws.Cells.Item(9, 6).Value2  = 452
$ws.Cells.Item(12, 6).Value2 = 784
$ws.Cells.Item(13, 6).Value2 = 45
$ws.Cells.Item(14, 6).Value2 = 87
$ws.Cells.Item(15, 6).Value2 = 347
$ws.Cells.Item(17, 6).Value2 = 301
$ws.Cells.Item(19, 6).Value2 = 400
$ws.Cells.Item(20, 6).Value2 = 272
$ws.Cells.Item(21, 6).Value2 = 1101
$ws.Cells.Item(23, 6).Value2 = 648
$ws.Cells.Item(24, 6).Value2 = 2239
$ws.Cells.Item(25, 6).Value2 = 760
$ws.Cells.Item(26, 6).Value2 = 58
$ws.Cells.Item(27, 6).Value2 = 561
$ws.Cells.Item(30, 6).Value2 = 558

# 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value2  = 292
$ws.Cells.Item(5, 6).Value2  = 315
$ws.Cells.Item(8, 6).Value2  = 30
$ws.Cells.Item(11, 6).Value2 = 4

# 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value2 = 474

# 全部类型 (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value2  = 474
$ws.Cells.Item(3, 6).Value2  = 7893
$ws.Cells.Item(6, 6).Value2  = 292
$ws.Cells.Item(7, 6).Value2  = 15247
$ws.Cells.Item(10, 6).Value2 = 654
$ws.Cells.Item(11, 6).Value2 = 452
$ws.Cells.Item(16, 6).Value2 = 315
$ws.Cells.Item(18, 6).Value2 = 784
$ws.Cells.Item(19, 6).Value2 = 45
$ws.Cells.Item(20, 6).Value2 = 87
$ws.Cells.Item(21, 6).Value2 = 347
$ws.Cells.Item(23, 6).Value2 = 21
$ws.Cells.Item(24, 6).Value2 = 30
$ws.Cells.Item(27, 6).Value2 = 301
$ws.Cells.Item(29, 6).Value2 = 400
$ws.Cells.Item(30, 6).Value2 = 273
$ws.Cells.Item(31, 6).Value2 = 1101
$ws.Cells.Item(33, 6).Value2 = 649
$ws.Cells.Item(34, 6).Value2 = 2239
$ws.Cells.Item(35, 6).Value2 = 761
$ws.Cells.Item(36, 6).Value2 = 58
$ws.Cells.Item(37, 6).Value2 = 561
$ws.Cells.Item(39, 6).Value2 = 4
$ws.Cells.Item(41, 6).Value2 = 558
